# Daily attendance processing - 2025-11-11 13:34:29
# Rotate the "Recorded By" (column G) comma-separated list so the last
# entry moves to the front, for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
